# Fixed opening and closing shift bug, adding serialization and deserialization of objects
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Opening shift: rows 12-17, columns B-F (Monday-Friday) should be "working" (1)
$ws.Range("B12:F13").Value = 1
$ws.Range("C14:C17").Value = 1
$ws.Range("E14:E17").Value = 1

# Closing shift: rows 23-25, columns B-F (Monday-Friday) should NOT be working (0)
$ws.Range("B23:F25").Value = 0

# Update the active cell selection to match the final cursor position
$ws.Range("F15").Select()
